# "Actualización automatica": add a "Varones" row (with its lookup URL) right
# before the existing "Mujeres" row, pushing "Mujeres" down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a copy of column B's current look (underlined hyperlink-blue text) in a
# scratch cell far outside the used range, so we can restore it later without
# picking up Excel's automatic "Hyperlink" named style.
$ws.Range("B1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# Insert a new row 4 (this shifts the old row 4, "Mujeres", down to row 5).
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Varones"
$ws.Range("B4").Value = "http://opendata.aragon.es/kos/iaest/sexo/varones"

# The row insert does not carry the hyperlink relationships down with it, so
# rebuild the hyperlinks for the whole column in the correct order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B1"), "http://opendata.aragon.es/kos/iaest/sexo/hombres")
$ws.Hyperlinks.Add($ws.Range("B2"), "http://opendata.aragon.es/kos/iaest/sexo/hombre")
$ws.Hyperlinks.Add($ws.Range("B3"), "http://opendata.aragon.es/kos/iaest/sexo/mujer")
$ws.Hyperlinks.Add($ws.Range("B4"), "http://opendata.aragon.es/kos/iaest/sexo/varones")
$ws.Hyperlinks.Add($ws.Range("B5"), "http://opendata.aragon.es/kos/iaest/sexo/mujeres")

# Hyperlinks.Add applies Excel's built-in "Hyperlink" cell style; drop it again
# and restore the plain look the sheet used before, so formatting is unchanged.
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("Z1").Copy()
$ws.Range("B1:B5").PasteSpecial(-4122)

$ws.Range("Z1").Clear()
